$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F4").Value = 832
$ws1.Range("F5").Value = 509
$ws1.Range("F7").Value = 9399
$ws1.Range("F10").Value = 683
$ws1.Range("F11").Value = 1866
$ws1.Range("F13").Value = 116
$ws1.Range("G13").Value = 59
$ws1.Range("F14").Value = 2531
$ws1.Range("F15").Value = 121
$ws1.Range("F16").Value = 3836
$ws1.Range("F17").Value = 291
$ws1.Range("F19").Value = 121
$ws1.Range("F20").Value = 201
$ws1.Range("F21").Value = 232
$ws1.Range("F22").Value = 13
$ws1.Range("F26").Value = 531
$ws1.Range("F27").Value = 2111
$ws1.Range("G27").Value = 29.9
$ws1.Range("F30").Value = 453
$ws1.Range("F31").Value = 4306
$ws1.Range("F33").Value = 109
$ws1.Range("F34").Value = 298
$ws1.Range("F35").Value = 97

# Sheet 2: 演出 (Performances)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("G2").Value = 0
$ws2.Range("F5").Value = 15

# Sheet 4: 全部类型 (All types)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("G5").Value = 0
$ws4.Range("F8").Value = 832
$ws4.Range("F9").Value = 509
$ws4.Range("F11").Value = 9399
$ws4.Range("F14").Value = 683
$ws4.Range("F15").Value = 1866
$ws4.Range("F17").Value = 116
$ws4.Range("G17").Value = 59
$ws4.Range("F19").Value = 2531
$ws4.Range("F20").Value = 121
$ws4.Range("F21").Value = 3836
$ws4.Range("F22").Value = 291
$ws4.Range("F24").Value = 121
$ws4.Range("F25").Value = 201
$ws4.Range("F26").Value = 232
$ws4.Range("F27").Value = 13
$ws4.Range("F28").Value = 15
$ws4.Range("F32").Value = 531
$ws4.Range("F33").Value = 2111
$ws4.Range("G33").Value = 29.9
$ws4.Range("F36").Value = 453
$ws4.Range("F37").Value = 4306
$ws4.Range("F39").Value = 109
$ws4.Range("F40").Value = 298
$ws4.Range("F41").Value = 97
